# Update for model 6 (row 7 in Sheet1 = "Model-6")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# P, C block (columns E:H) -> Acctr, Sentr, Spectr, Mcctr
$ws.Range("E7:G7").NumberFormat = "0.00"
$ws.Range("E7:G7").Value = 100
$ws.Range("H7").NumberFormat = "0.00"
$ws.Range("H7").Value = 1

# P, C block (columns J:M) -> Acccv, Sencv, Speccv, Mcccv
$ws.Range("J7:L7").NumberFormat = "0.00"
$ws.Range("J7:L7").Value = 100
$ws.Range("M7").NumberFormat = "0.00"
$ws.Range("M7").Value = 1

# P, C block (columns O:R) -> Accext, Senext, Specext, Mccext
$ws.Range("O7:Q7").NumberFormat = "0.00"
$ws.Range("O7:Q7").Value = 100
$ws.Range("R7").NumberFormat = "0.00"
$ws.Range("R7").Value = 1

# Match the author's final selection after entering the new row
$null = $ws.Range("O7:R7").Select()
